$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.969.11"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "3.510.70"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.89"
$ws.Range("E5").Value = "  -3.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.97"
$ws.Range("E6").Value = "  -5.55%  "
$ws.Range("D7").Value = "3.506.88"
$ws.Range("E7").Value = "  -3.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.59"
$ws.Range("E11").Value = "  +4.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  -4.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.13"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "4.109.97"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").Value = "3.517.96"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").Value = "68.134.54"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.46"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.96"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.23"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.628"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "3.655.46"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("E27").Value = "  -8.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.73"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  -6.49%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.66"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.67"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.21"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.508.24"
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.85"
$ws.Range("E37").Value = "  -6.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.04"
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.66"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0908"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.23"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.899"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "47.03"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.52"
$ws.Range("E49").Value = "  -10.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.64"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -2.08%  "
